$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column C (Länge), shifting Länge..Nominatim Wichtigkeit
# one column to the right (to D..G), to make room for the new "Webseite" column.
$ws.Columns("C:C").Insert()

# Give the new header cell (C1) the same formatting as the other header cells
# (bold header style used across row 1), then set its text.
$ws.Range("A1").Copy() | Out-Null
$ws.Range("C1").PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$ws.Range("C1").Value = "Webseite"

# Fill in the website addresses for each location (rows 2-7).
$ws.Range("C2").Value = "https://www.allerweltshaus.de/"
$ws.Range("C3").Value = "https://solawi-alfter.de/"
$ws.Range("C4").Value = "https://alsenstrasse.com/der-alsengarten/"
$ws.Range("C5").Value = "https://urbane-gaerten.de/urbane-gaerten/gaerten-im-ueberblick/demogarten,-k%C3%B6ln"
$ws.Range("C6").Value = "https://www.urbanes-gaertnern-freiburg.de/de/organisations/fridas-klima-garten"
$ws.Range("C7").Value = "https://www.weltacker-berlin.de/"

# Match the new column's width to its neighbor (column B).
$ws.Columns("C:C").ColumnWidth = $ws.Columns("B:B").ColumnWidth

# Update the active selection to match the saved workbook state.
$ws.Range("C7").Select()

$excel.CutCopyMode = $false
